# Update "想去人数" (column F) figures across the workbook's sheets to
# reflect newly scraped counts, as published by the gh-pages build at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1470
$ws1.Range("F3").Value = 1443
$ws1.Range("F6").Value = 705
$ws1.Range("F8").Value = 634
$ws1.Range("F12").Value = 33287
$ws1.Range("F13").Value = 7083
$ws1.Range("F14").Value = 114
$ws1.Range("F15").Value = 368
$ws1.Range("F16").Value = 577
$ws1.Range("F17").Value = 444
$ws1.Range("F20").Value = 19
$ws1.Range("F21").Value = 50
$ws1.Range("F22").Value = 450
$ws1.Range("F24").Value = 805
$ws1.Range("F27").Value = 390
$ws1.Range("F30").Value = 210
$ws1.Range("F31").Value = 51
$ws1.Range("F34").Value = 133
$ws1.Range("F35").Value = 742
$ws1.Range("F38").Value = 793

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1200
$ws2.Range("F20").Value = 13

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1459
$ws3.Range("F3").Value = 358

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1459
$ws4.Range("F3").Value = 358
$ws4.Range("F4").Value = 1200
$ws4.Range("F5").Value = 1470
$ws4.Range("F7").Value = 1443
$ws4.Range("F9").Value = 705
$ws4.Range("F11").Value = 634
$ws4.Range("F20").Value = 7083
$ws4.Range("F21").Value = 368
$ws4.Range("F23").Value = 577
$ws4.Range("F24").Value = 444
$ws4.Range("F27").Value = 19
$ws4.Range("F29").Value = 50
$ws4.Range("F31").Value = 450
$ws4.Range("F33").Value = 806
$ws4.Range("F36").Value = 390
$ws4.Range("F39").Value = 210
$ws4.Range("F40").Value = 51
$ws4.Range("F44").Value = 133
$ws4.Range("F45").Value = 793
$ws4.Range("F48").Value = 13
